$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-24"

# Update header label for the "2022" column (shared string used by I1)
$ws.Range("I1").Value = "2022 (through 07-24)"

# Update August value for 2022 column (row 8 = August)
$ws.Range("I8").Value = 136

# Update Total value for 2022 column (row 14 = Total)
$ws.Range("I14").Value = 942
